# Global_Material_Data.xlsx - "Change Material Global only SAP"
#
# The existing test case in row 6 (Create_Material_with_Questionnaire_only_Global_SAP)
# is switched off (Execute = N) and a new test case is appended in row 7 for
# "3.Change_Material_global_only" (Execute = Y), reusing the same row layout/styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stop executing the existing "Create_Material_with_Questionnaire_only_Global_SAP" case
$ws.Range("C6").Value = "N"

# Duplicate row 6 into row 7 so the new test case keeps the same formatting/styles
$ws.Range("A6:AK6").Copy($ws.Range("A7:AK7"))

# Fill in the new "Change Material Global only" test case details in row 7
$ws.Range("B7").Value = "3.Change_Material_global_only"
$ws.Range("C7").Value = "Y"
$ws.Range("E7").Value = "Changing Material Descript"
$ws.Range("G7").Value = "CMG0012"
$ws.Range("K7").Value = "KG"
$ws.Range("N7").Value = "KG"

# Leave the selection where the author left it after editing
$ws.Range("M14").Select()
